$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand/receptor TPM-derived metrics for Egf-Egfr LR pairs (rows 2-10)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.041452
$ws.Range("H2").Value = 0.124356
$ws.Range("I2").Value = 0.05439747478414846
$ws.Range("J2").Value = 0.05439747478414846
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.017772254836
$ws.Range("R2").Value = 0.159950293524
$ws.Range("S2").Value = 0.000224145009643676
$ws.Range("T2").Value = 0.000224145009643676
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.041452
$ws.Range("H3").Value = 0.124356
$ws.Range("I3").Value = 0.05439747478414846
$ws.Range("J3").Value = 0.05439747478414846
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 3.325537976124
$ws.Range("R3").Value = 29.929841785116
$ws.Range("S3").Value = 0.04194193413313065
$ws.Range("T3").Value = 0.04194193413313064
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.041452
$ws.Range("H4").Value = 0.124356
$ws.Range("I4").Value = 0.05439747478414846
$ws.Range("J4").Value = 0.05439747478414846
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 0.9698162840386668
$ws.Range("R4").Value = 8.728346556348001
$ws.Range("S4").Value = 0.01223139564137414
$ws.Range("T4").Value = 0.01223139564137414
$ws.Range("I5").Value = 0.5204718857143857
$ws.Range("J5").Value = 0.5204718857143856
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 0.1700439041443333
$ws.Range("R5").Value = 1.530395137299
$ws.Range("S5").Value = 0.002144606460237903
$ws.Range("T5").Value = 0.002144606460237903
$ws.Range("I6").Value = 0.5204718857143857
$ws.Range("J6").Value = 0.5204718857143856
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 31.81855459864899
$ws.Range("S6").Value = 0.4012979947212597
$ws.Range("T6").Value = 0.4012979947212596
$ws.Range("I7").Value = 0.5204718857143857
$ws.Range("J7").Value = 0.5204718857143856
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 9.279145992585891
$ws.Range("R7").Value = 83.51231393327301
$ws.Range("S7").Value = 0.1170292845328882
$ws.Range("T7").Value = 0.1170292845328881
$ws.Range("G8").Value = 0.3239583333333333
$ws.Range("H8").Value = 0.971875
$ws.Range("I8").Value = 0.4251306395014658
$ws.Range("J8").Value = 0.4251306395014658
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 0.1388948677083333
$ws.Range("R8").Value = 1.250053809375
$ws.Range("S8").Value = 0.001751752478750101
$ws.Range("T8").Value = 0.001751752478750101
$ws.Range("G9").Value = 0.3239583333333333
$ws.Range("H9").Value = 0.971875
$ws.Range("I9").Value = 0.4251306395014658
$ws.Range("J9").Value = 0.4251306395014658
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 25.989958028125
$ws.Range("R9").Value = 233.909622253125
$ws.Range("S9").Value = 0.3277872980446166
$ws.Range("T9").Value = 0.3277872980446166
$ws.Range("G10").Value = 0.3239583333333333
$ws.Range("H10").Value = 0.971875
$ws.Range("I10").Value = 0.4251306395014658
$ws.Range("J10").Value = 0.4251306395014658
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 7.579370525347223
$ws.Range("R10").Value = 68.21433472812501
$ws.Range("S10").Value = 0.09559158897809913
$ws.Range("T10").Value = 0.09559158897809912
